$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per the commit diff (Price/Volume(1h) refresh for Fri May 17
# 2024 run, plus three row pairs whose coins swapped rank/row position).

$ws.Range("D2").Value = '65.407.92'
$ws.Range("E2").Value = '  -0.76%  '
$ws.Range("D3").Value = '2.950.19'
$ws.Range("E3").Value = '  -1.79%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.93'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.97%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.517'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = '2.946.89'
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.71'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.30%  '
$ws.Range("E11").Value = '  -3.89%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("E13").Value = '  -2.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.48'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("D16").Value = '65.401.76'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").Value = '3.439.80'
$ws.Range("E17").Value = '  -1.85%  '
$ws.Range("D19").Value = '2.950.05'
$ws.Range("E19").Value = '  -1.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +13.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '444.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.695'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("E25").Value = '  -2.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.21'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.22%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.73%  '
$ws.Range("B29").Value = 'NEARProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.29%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0000102'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.46%  '
$ws.Range("E33").Value = '  +2.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.15'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.84%  '
$ws.Range("B37").Value = 'Mantle'
$ws.Range("C37").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.969'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.09'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '44.85'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.54%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.121'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.54%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.96'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.22%  '
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("E43").Value = '  -4.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.53'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '383.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("E46").Value = '  -1.45%  '
$ws.Range("D47").Value = '2.678.09'
$ws.Range("E47").Value = '  -4.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.32'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.34%  '
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.17'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.80%  '
